# Update "想去人数" (column F) figures on sheets "展览" and "全部类型"
# to reflect freshly scraped counts (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Map of sheet name -> { row -> new F value }
$updates = @{
    "展览" = @{
        2  = 138
        4  = 19
        5  = 6758
        6  = 86
        7  = 6
        8  = 436
        10 = 6317
        13 = 1275
        14 = 11
        15 = 103
        16 = 389
        18 = 18
        19 = 372
        20 = 45
        21 = 7
        22 = 4674
        23 = 67
        24 = 48
        25 = 127
        26 = 196
        27 = 88
    }
    "全部类型" = @{
        2  = 138
        4  = 19
        5  = 6758
        6  = 86
        7  = 6
        8  = 436
        10 = 6317
        13 = 1275
        14 = 11
        15 = 103
        16 = 389
        18 = 18
        19 = 372
        20 = 45
        21 = 7
        22 = 4674
        24 = 67
        25 = 48
        26 = 127
        27 = 196
        28 = 88
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
